$d = $word.ActiveDocument

# --- Programa paragraph: split the big run into 5 pieces using manual line breaks ---
$d.Content.Find.Execute("expansão linear.2. Estrutura cristalina", $false, $false, $false, $false, $false, $true, 1, $false, "expansão linear.^l2. Estrutura cristalina", 2) | Out-Null
$d.Content.Find.Execute("difração de raios-X;3. Defeitos em cristais", $false, $false, $false, $false, $false, $true, 1, $false, "difração de raios-X;^l3. Defeitos em cristais", 2) | Out-Null
$d.Content.Find.Execute("inclusões).4. Relação microestrutura", $false, $false, $false, $false, $false, $true, 1, $false, "inclusões).^l4. Relação microestrutura", 2) | Out-Null
$d.Content.Find.Execute("magnéticas).Em todos os itens", $false, $false, $false, $false, $false, $true, 1, $false, "magnéticas).^lEm todos os itens", 2) | Out-Null

# --- Norma de recuperação paragraph: split formula sentence before "MR = " ---
$d.Content.Find.Execute("calculada pela fórmula: MR = (NF + PR) / 2", $false, $false, $false, $false, $false, $true, 1, $false, "calculada pela fórmula: ^lMR = (NF + PR) / 2", 2) | Out-Null

# --- Bibliografia paragraph: split each numbered reference onto its own line ---
$d.Content.Find.Execute("LTC Editora, 2013.2. ASKELAND", $false, $false, $false, $false, $false, $true, 1, $false, "LTC Editora, 2013.^l2. ASKELAND", 2) | Out-Null
$d.Content.Find.Execute("São Paulo, 2008.3. SHACKELFORD", $false, $false, $false, $false, $false, $true, 1, $false, "São Paulo, 2008.^l3. SHACKELFORD", 2) | Out-Null
$d.Content.Find.Execute("6a. ed., Pearson, 2008.4. PADILHA, A.F., Materiais", $false, $false, $false, $false, $false, $true, 1, $false, "6a. ed., Pearson, 2008.^l4. PADILHA, A.F., Materiais", 2) | Out-Null
$d.Content.Find.Execute("Hemus Editora, 1997.5. PADILHA, A.F., Técnicas", $false, $false, $false, $false, $false, $true, 1, $false, "Hemus Editora, 1997.^l5. PADILHA, A.F., Técnicas", 2) | Out-Null
$d.Content.Find.Execute("Ed. Hemus, 1985.6. REED-HILL", $false, $false, $false, $false, $false, $true, 1, $false, "Ed. Hemus, 1985.^l6. REED-HILL", 2) | Out-Null
$d.Content.Find.Execute("Guanabara Dois, 1982.7. BRANDON", $false, $false, $false, $false, $false, $true, 1, $false, "Guanabara Dois, 1982.^l7. BRANDON", 2) | Out-Null
$d.Content.Find.Execute("1st. ed., Wiley, 1999.8. ASHBY, M.F., JONES", $false, $false, $false, $false, $false, $true, 1, $false, "1st. ed., Wiley, 1999.^l8. ASHBY, M.F., JONES", 2) | Out-Null
$d.Content.Find.Execute("Elsevier Editora, 2007.9. ASHBY, M.F., SHERCLIFF", $false, $false, $false, $false, $false, $true, 1, $false, "Elsevier Editora, 2007.^l9. ASHBY, M.F., SHERCLIFF", 2) | Out-Null
